$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The workbook originally carried a duplicate "Text" style (numFmtId 49) as
# both style index 1 and style index 3; re-applying the format explicitly
# collapses rows that referenced the duplicate (B38:E39, D40:E40) onto the
# single surviving style, matching the de-duplicated style table.
$ws.Range("B38:E40").NumberFormat = "@"

# New rows 41-42: same column formatting as the rest of the table
# (Text for B:E, #,##0 for F:K).
$ws.Range("B41:E42").NumberFormat = "@"
$ws.Range("F41:K42").NumberFormat = "#,##0"

# Row 41 - RMO No. 11-2024
$ws.Cells.Item(41,1).Value2  = "https://bir-cdn.bir.gov.ph/BIR/pdf/RMO%20No.%2011-2024%20Annexes.pdf"
$ws.Cells.Item(41,2).Value2  = "RMO No. 11-2024"
$ws.Cells.Item(41,3).Value2  = "March 14, 2024"
$ws.Cells.Item(41,4).Value2  = "2024"
$ws.Cells.Item(41,5).Value2  = "Initial"
$ws.Cells.Item(41,6).Value2  = 3055169
$ws.Cells.Item(41,7).Formula = "=1649267+72127"
$ws.Cells.Item(41,8).Value2  = 326195
$ws.Cells.Item(41,9).Value2  = 599235
$ws.Cells.Item(41,10).Value2 = 163216
$ws.Cells.Item(41,11).Formula = "=229242+15887"

# Row 42 - RMO No. 29-2024
$ws.Cells.Item(42,1).Value2  = "https://bir-cdn.bir.gov.ph/BIR/pdf/CY2024%20Goal%20_RMO%2029-2024%20Annexes.xlsx"
$ws.Cells.Item(42,2).Value2  = "RMO No. 29-2024"
$ws.Cells.Item(42,3).Value2  = "July 22, 2024"
$ws.Cells.Item(42,4).Value2  = "2024"
$ws.Cells.Item(42,5).Value2  = "Revision 1"
$ws.Cells.Item(42,6).Value2  = 3046751.2749999999
$ws.Cells.Item(42,7).Formula = "=1467398.307+ 85227.022"
$ws.Cells.Item(42,8).Value2  = 457812.31400000001
$ws.Cells.Item(42,9).Value2  = 620708.23899999994
$ws.Cells.Item(42,10).Value2 = 144379.228
$ws.Cells.Item(42,11).Formula = "=255339.165+ 15887"
